# atualizei dados bibi e add
# Update faturamento_diario sheet:
#  - correct B2 (Jul day1) and B6 (Jul day7) totals
#  - insert a new daily record (Jul day8) as row 7, shifting existing rows down
#  - correct the June day30 total (now row 28 after the shift)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing July values
$ws.Range("B2").Value = 17972.28
$ws.Range("B6").Value = 15070.96

# Insert a new row at position 7 (pushes rows 7..68 down to 8..69)
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new daily record (Jul day 8)
$ws.Range("A7").Value = 8
$ws.Range("B7").Value = 19918.15
$ws.Range("C7").Value = 7
$ws.Range("D7").Value = 2025
$ws.Range("E7").Value = "07/2025"

# Correct the June day-30 total, which is now on row 28 after the insert/shift
$ws.Range("B28").Value = 111900.66
